$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Edits inside the existing "مرداد 99" block (rows 104-112)
# ---------------------------------------------------------------------------

# Row 107: the (previously empty, styled) E107 cell disappears entirely.
$ws.Range("E107").Clear()

# Row 108: hours worked on "* Online Tracking" go from 3 to 5.
$ws.Range("C108").Value = 5

# Row 109: "* Other" task becomes "* Meetings & other" (hours unchanged, 5).
$ws.Range("B109").Value = "* Meetings & other"

# Rows 110/112 (C110, E110, D112) are formulas and recompute automatically.

# ---------------------------------------------------------------------------
# 2) Append a brand-new monthly block (شهریور و مهر 99) in rows 113-123,
#    built by duplicating the previous block's formatting (rows 102-112)
#    and then overwriting its contents.
# ---------------------------------------------------------------------------

$src = $ws.Range("A102:G112")
$dst = $ws.Range("A113")
$src.Copy($dst)

# Row 116 (set the new bullet string first so the shared-string table is
# rebuilt in the same order as the source edit)
$ws.Range("B116").Value = "* Code Refactoring"
$ws.Range("C116").Value = 1
$ws.Range("E116").Value = "• "

# Row 117
$ws.Range("B117").Value = "* GUI"
$ws.Range("C117").Value = 2
$ws.Range("E117").Value = "• "

# Row 115 - block header: month title + column captions
$ws.Range("A115").Value = "شهریور و مهر 99"
$ws.Range("B115").Value = "Activity"
$ws.Range("C115").Value = "Hours"
$ws.Range("E115").Value = "Tasks Done"

# Row 118 (no Tasks-Done entry this row)
$ws.Range("B118").Value = "* Registration"
$ws.Range("C118").ClearContents()
$ws.Range("E118").Clear()

# Row 119
$ws.Range("B119").Value = "* Online Tracking"
$ws.Range("C119").ClearContents()
$ws.Range("E119").Clear()

# Row 120
$ws.Range("B120").Value = "* Meetings & other"
$ws.Range("C120").ClearContents()
$ws.Range("E120").ClearContents()

# Row 121 - totals for the new block
# (Clear() first: the freshly pasted formula cells keep a stale cached
#  result tied to the copied-from range until the cell is reset, even
#  though `.Formula` text reads back correctly.)
$ws.Range("B121").Value = "• Total Hours"
$ws.Range("C121").Clear()
$ws.Range("C121").Formula = "=SUM(C116:C120)"
$ws.Range("E121").Clear()
$ws.Range("E121").Formula = "=C121*40000"

# Row 122 - Paid
$ws.Range("C122").Value = "@Parsiss"
$ws.Range("D122").Value = 0

# Row 123 - Not Paid
$ws.Range("C123").Value = "@Home"
$ws.Range("D123").Clear()
$ws.Range("D123").Formula = "=C121-D122"

# ---------------------------------------------------------------------------
# 3) Update the selected cell shown when the workbook is reopened.
# ---------------------------------------------------------------------------
$ws.Range("E119").Select()
